# Updated symbol list on Sat Jan 14 03:27:59 UTC 2023 with GitHub Actions
#
# Writes new Price (column D) / Volume(1h) (column E) figures, and swaps the
# FTXToken / GateToken rows (7 and 8) to reflect the refreshed ranking order.
# All of these sheet cells are plain text (t="inlineStr" in the original
# OOXML) even though many of them look numeric ("305.17", "5.60%", ...), so
# each write forces the cell to Text format first (and resets the style back
# to Normal afterwards) to stop Excel's automatic number/percentage parsing
# from turning the literal string into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Addr,
        [string]$Text
    )
    $rng = $ws.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2 - BNB
Set-TextCell "D2" "306.11"
Set-TextCell "E2" "6.16%"

# Row 3 - OKB
Set-TextCell "D3" "32.50"
Set-TextCell "E3" "10.69%"

# Row 4 - HuobiToken
Set-TextCell "D4" "5.316"
Set-TextCell "E4" "3.49%"

# Row 5 - Cronos
Set-TextCell "D5" "0.07417"
Set-TextCell "E5" "11.45%"

# Row 6 - KuCoinToken
Set-TextCell "D6" "7.746"
Set-TextCell "E6" "5.52%"

# Row 7 - was FTXToken, now GateToken
Set-TextCell "B7" "GateToken"
Set-TextCell "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D7" "3.715"
Set-TextCell "E7" "9.05%"

# Row 8 - was GateToken, now FTXToken
Set-TextCell "B8" "FTXToken"
Set-TextCell "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D8" "1.582"
Set-TextCell "E8" "16.79%"

# Row 9 - MXToken
Set-TextCell "D9" "0.9219"
Set-TextCell "E9" "0.40%"

# Row 10 - One
Set-TextCell "D10" "0.01623"
Set-TextCell "E10" "2,412.97%"

# Row 11 - WazirX
Set-TextCell "D11" "0.1668"
Set-TextCell "E11" "6.12%"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextCell "D12" "0.07416"
Set-TextCell "E12" "13.28%"

# Row 13 - MandalaExchangeToken
Set-TextCell "D13" "0.07978"
Set-TextCell "E13" "3.84%"

# Row 14 - BitrueCoin
Set-TextCell "D14" "0.03106"
Set-TextCell "E14" "6.85%"

# Row 15 - BitMartToken
Set-TextCell "D15" "0.09812"
Set-TextCell "E15" "9.10%"

# Row 16 - BitForexToken
Set-TextCell "D16" "0.001521"
Set-TextCell "E16" "-3.67%"

# Row 17 - CoinExToken (price unchanged, only Volume(1h))
Set-TextCell "E17" "1.83%"

# Row 18 - TigerCash
Set-TextCell "D18" "0.006155"
Set-TextCell "E18" "-1.52%"

# Row 19 - LEO
Set-TextCell "D19" "3.470"
Set-TextCell "E19" "0.38%"

# Row 20 - BTSEToken
Set-TextCell "D20" "2.240"
Set-TextCell "E20" "0.56%"

# Row 21 - BitpandaEcosystemToken
Set-TextCell "D21" "0.3272"
Set-TextCell "E21" "1.85%"

# Row 22 - ProBitToken
Set-TextCell "D22" "0.1312"
Set-TextCell "E22" "0.23%"

# Row 23 - MCDex
Set-TextCell "D23" "4.255"
Set-TextCell "E23" "4.79%"

# Row 24 - ZBToken
Set-TextCell "D24" "0.1639"
Set-TextCell "E24" "5.75%"

# Row 25 - BitKan (price unchanged, only Volume(1h))
Set-TextCell "E25" "3.17%"

# Row 26 - HotbitToken
Set-TextCell "D26" "0.004536"
Set-TextCell "E26" "9.94%"

# Row 27 - NitroEx
Set-TextCell "D27" "0.0001169"
Set-TextCell "E27" "-6.40%"

# Row 28 - UpBots
Set-TextCell "D28" "0.0001666"
Set-TextCell "E28" "3.03%"

# Row 40 - IDEX
Set-TextCell "D40" "0.04510"
Set-TextCell "E40" "7.33%"

# Row 41 - KickToken
Set-TextCell "D41" "0.007308"
Set-TextCell "E41" "8.19%"

# Row 42 - BKEXToken
Set-TextCell "D42" "0.1366"
Set-TextCell "E42" "10.27%"

# Row 43 - CEJI
Set-TextCell "D43" "0.002178"
Set-TextCell "E43" "10.10%"

# Row 44 - LocalTraders
Set-TextCell "D44" "0.01374"
Set-TextCell "E44" "8.00%"

# Row 45 - CoinLion
Set-TextCell "D45" "0.00006009"
Set-TextCell "E45" "6.33%"

# Row 47 - CoinbaseStockToken
Set-TextCell "D47" "0.01300"
Set-TextCell "E47" "-0.46%"
